$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New data row appended below the existing table (row 96 -> row 97).
# Column A holds a date-like text label ("2025/10/13"); prefix with an
# apostrophe so it is stored as literal text (matching the workbook's
# existing rows) instead of being auto-converted to a date serial, then
# reset the style so no extra "text" formatting is stamped on the cell.
$ws.Range("A97").Value = "'2025/10/13"
$ws.Range("A97").Style = "Normal"

$ws.Range("B97").Value = "月"
$ws.Range("C97").Value = 0
$ws.Range("D97").Value = 201
